$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "verbrauch" label in M2 becomes "leds"
$ws.Range("M2").Value = "leds"

# Row 3: change existing anzahl (L3) value, add a new "watt/leds" count (M3)
# and a total (N3 = M3*L3)
$ws.Range("L3").Value = 22
$ws.Range("M3").Value = 24
$ws.Range("N3").Formula = "=M3*L3"

# Row 4: add L4/M4 counts and N4 = M4*L4
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 35
$ws.Range("N4").Formula = "=M4*L4"

# Row 5: add a new label in K5 ("matric"), L5/M5 counts, N5 = M5*L5
$ws.Range("K5").Value = "matric"
$ws.Range("L5").Value = 8
$ws.Range("M5").Value = 30
$ws.Range("N5").Formula = "=M5*L5"

# Row 6: totals - N6 sums N3:N5, O6 is a percentage factor, P6 = O6*N6
$ws.Range("N6").Formula = "=N3+N4+N5"
$ws.Range("O6").Value = 0.03
$ws.Range("P6").Formula = "=O6*N6"

# Update the current selection to match the saved view state
$ws.Range("M6").Select()
